# "changed data file and added formatting"
#
# The "Renewal capex" sheet had three extra metadata columns
# (B: "Available from", C: "Available until", D: "Technology type")
# inserted before the yearly data columns. The updated data file drops
# those three columns so that "Renewal capex" has the same layout as the
# other two sheets (label in column A, then the year columns 2020-2050
# in B:AF).
$wb = $excel.ActiveWorkbook

$renewalCapex = $wb.Worksheets.Item("Renewal capex")
$otherOpex    = $wb.Worksheets.Item("Other Opex")

# Remove the "Available from" / "Available until" / "Technology type"
# columns - the remaining year columns shift left into B:AF.
$renewalCapex.Range("B1:D1").EntireColumn.Delete()

# Restore/update the selections on the affected sheets (cosmetic, but part
# of the saved view state after the edit).
$otherOpex.Activate()
$otherOpex.Range("F15").Select()

$renewalCapex.Activate()
$renewalCapex.Range("C2").Select()
